$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E, shifting boson/value/stat_u/syst_u/obs/diff/target/col right by one.
$ws.Columns("E:E").Insert()

# Header for new column
$ws.Range("E1").Value = "pt_max"

# Fill pt_max = 50 for each data row (rows 2-9)
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 5).Value = 50
}

# Update selection to match target state
$ws.Range("E15").Select()
